$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.681.86"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.628.77"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.47"
$ws.Range("E5").Value = "  -0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.83"
$ws.Range("E6").Value = "  -1.18%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -3.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.625.94"
$ws.Range("E9").Value = "  +0.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.118"
$ws.Range("E10").Value = "  -4.00%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.81"
$ws.Range("E11").Value = "  +0.16%  "

$ws.Range("E12").Value = "  -2.26%  "

$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.38"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.102.07"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("E16").Value = "  -2.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.583.54"
$ws.Range("E17").Value = "  -1.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.633.71"
$ws.Range("E18").Value = "  +0.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.12"
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.59"
$ws.Range("E20").Value = "  +2.92%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.53"
$ws.Range("E21").Value = "  -3.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "344.80"
$ws.Range("E22").Value = "  -0.36%  "

$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.96"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.89"
$ws.Range("E25").Value = "  +9.27%  "

$ws.Range("E26").Value = "  -3.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "606.62"
$ws.Range("E27").Value = "  +8.90%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.24"
$ws.Range("E28").Value = "  -1.65%  "

$ws.Range("E29").Value = "  +2.50%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.97"
$ws.Range("E30").Value = "  -0.37%  "

$ws.Range("E31").Value = "  -0.81%  "

$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("E33").Value = "  -1.15%  "

$ws.Range("E34").Value = "  +0.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.62"
$ws.Range("E35").Value = "  +2.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.40"
$ws.Range("E36").Value = "  +0.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.402"
$ws.Range("E37").Value = "  -2.40%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.72"
$ws.Range("E39").Value = "  -1.66%  "

$ws.Range("E40").Value = "  -2.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.07"
$ws.Range("E41").Value = "  -1.02%  "

$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.53"
$ws.Range("E43").Value = "  +2.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.70"
$ws.Range("E44").Value = "  -0.64%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.81"
$ws.Range("E45").Value = "  +8.15%  "

$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "158.86"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("E47").Value = "  -2.37%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0586"
$ws.Range("E48").Value = "  -2.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.630"
$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0997"
$ws.Range("E50").Value = "  -1.79%  "

$ws.Range("E51").Value = "  -0.59%  "
